$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed cryptos list snapshot.
# Column D prices that are plain decimals (no thousands separators) get
# auto-converted to numbers by Excel on assignment, so we force those
# specific cells to Text format first to preserve the original string look.

$ws.Range("D2").Value = '69.114.37'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '3.807.11'
$ws.Range("E3").Value = '  +1.64%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.55'
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.76'
$ws.Range("E6").Value = '  -3.06%  '

$ws.Range("D7").Value = '3.804.85'
$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("E10").Value = '  +1.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.30'
$ws.Range("E11").Value = '  -1.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("E12").Value = '  -0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.31'
$ws.Range("E13").Value = '  -2.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000246'
$ws.Range("E14").Value = '  -1.30%  '

$ws.Range("D15").Value = '4.446.57'
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("D16").Value = '3.803.73'
$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("D17").Value = '69.221.01'
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.37'
$ws.Range("E19").Value = '  +1.67%  '

$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.114'
$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.46'
$ws.Range("E21").Value = '  +4.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.01'
$ws.Range("E22").Value = '  -1.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("E23").Value = '  -1.04%  '

$ws.Range("E24").Value = '  +3.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.82'
$ws.Range("E25").Value = '  -0.58%  '

$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("E27").Value = '  -1.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  -2.25%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.39'
$ws.Range("E32").Value = '  -4.46%  '

$ws.Range("D33").Value = '3.957.52'
$ws.Range("E33").Value = '  +1.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.79'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("D35").Value = '3.752.59'
$ws.Range("E35").Value = '  +2.01%  '

$ws.Range("E36").Value = '  -1.47%  '

$ws.Range("E37").Value = '  +5.35%  '

$ws.Range("E38").Value = '  +0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.05'
$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.320'
$ws.Range("E42").Value = '  -1.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '429.62'
$ws.Range("E43").Value = '  -1.51%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.99'
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.56'
$ws.Range("E45").Value = '  -0.16%  '

$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.38'
$ws.Range("E47").Value = '  -1.31%  '

$ws.Range("D48").Value = '2.835.74'
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.49'
$ws.Range("E49").Value = '  +0.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.57'
$ws.Range("E50").Value = '  -2.99%  '

$ws.Range("E51").Value = '  -0.92%  '
